# Apply crypto price/volume updates from the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text entry so numeric-looking strings (e.g. "405.60",
# "68.088.61") are not auto-converted to numbers and keep their exact formatting,
# matching how the source data is stored (inline/shared text strings).

$ws.Range("D2").Value = "'68.088.61"
$ws.Range("E2").Value = "'  +3.90%  "
$ws.Range("D3").Value = "'3.612.80"
$ws.Range("E3").Value = "'  +4.16%  "
$ws.Range("E4").Value = "'  -0.20%  "
$ws.Range("D5").Value = "'202.13"
$ws.Range("E5").Value = "'  +10.87%  "
$ws.Range("D6").Value = "'574.52"
$ws.Range("E6").Value = "'  +2.46%  "
$ws.Range("D7").Value = "'0.617"
$ws.Range("E7").Value = "'  +2.93%  "
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("D9").Value = "'0.683"
$ws.Range("E9").Value = "'  +5.61%  "
$ws.Range("D10").Value = "'60.54"
$ws.Range("E10").Value = "'  +17.55%  "
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "'  +5.55%  "
$ws.Range("D12").Value = "'0.0000283"
$ws.Range("E12").Value = "'  +13.44%  "
$ws.Range("D13").Value = "'10.39"
$ws.Range("E13").Value = "'  +9.71%  "
$ws.Range("D14").Value = "'4.201.66"
$ws.Range("E14").Value = "'  +4.04%  "
$ws.Range("D15").Value = "'3.634.34"
$ws.Range("E15").Value = "'  +4.64%  "
$ws.Range("B16").Value = "'TRON"
$ws.Range("C16").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.127"
$ws.Range("E16").Value = "'  +1.67%  "
$ws.Range("B17").Value = "'Chainlink"
$ws.Range("C17").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'19.28"
$ws.Range("E17").Value = "'  +9.57%  "
$ws.Range("D18").Value = "'67.988.85"
$ws.Range("E18").Value = "'  +3.91%  "
$ws.Range("D19").Value = "'12.31"
$ws.Range("E19").Value = "'  +5.49%  "
$ws.Range("D20").Value = "'1.07"
$ws.Range("E20").Value = "'  +3.49%  "
$ws.Range("D21").Value = "'405.60"
$ws.Range("E21").Value = "'  +7.17%  "
$ws.Range("D22").Value = "'12.97"
$ws.Range("E22").Value = "'  +22.17%  "
$ws.Range("D23").Value = "'4.21"
$ws.Range("E23").Value = "'  +2.82%  "
$ws.Range("D24").Value = "'85.24"
$ws.Range("E24").Value = "'  +2.48%  "
$ws.Range("D25").Value = "'3.96"
$ws.Range("E25").Value = "'  +14.72%  "
$ws.Range("D26").Value = "'2.91"
$ws.Range("E26").Value = "'  +3.62%  "
$ws.Range("D27").Value = "'12.56"
$ws.Range("E27").Value = "'  +5.62%  "
$ws.Range("D28").Value = "'6.12"
$ws.Range("E28").Value = "'  +2.20%  "
$ws.Range("D29").Value = "'9.31"
$ws.Range("E29").Value = "'  +8.17%  "
$ws.Range("D30").Value = "'7.73"
$ws.Range("E30").Value = "'  +6.74%  "
$ws.Range("D31").Value = "'31.59"
$ws.Range("E31").Value = "'  +4.37%  "
$ws.Range("D32").Value = "'679.06"
$ws.Range("E32").Value = "'  +11.98%  "
$ws.Range("D33").Value = "'12.20"
$ws.Range("E33").Value = "'  +2.83%  "
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "'  +3.64%  "
$ws.Range("D35").Value = "'63.64"
$ws.Range("E35").Value = "'  +2.28%  "
$ws.Range("D36").Value = "'41.68"
$ws.Range("E36").Value = "'  +2.94%  "
$ws.Range("D37").Value = "'0.412"
$ws.Range("E37").Value = "'  +5.10%  "
$ws.Range("E38").Value = "'  -0.22%  "
$ws.Range("D39").Value = "'0.0₃0761"
$ws.Range("E39").Value = "'  +7.55%  "
$ws.Range("D40").Value = "'3.19"
$ws.Range("E40").Value = "'  +16.69%  "
$ws.Range("B41").Value = "'Maker"
$ws.Range("C41").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'3.194.96"
$ws.Range("E41").Value = "'  +8.90%  "
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.134"
$ws.Range("E42").Value = "'  +4.91%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "'  -0.13%  "
$ws.Range("D44").Value = "'2.70"
$ws.Range("E44").Value = "'  +10.88%  "
$ws.Range("B45").Value = "'WEMIXToken"
$ws.Range("C45").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.84"
$ws.Range("E45").Value = "'  +15.79%  "
$ws.Range("B46").Value = "'dogwifhat"
$ws.Range("C46").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.83"
$ws.Range("E46").Value = "'  +24.92%  "
$ws.Range("D47").Value = "'0.0413"
$ws.Range("E47").Value = "'  +4.77%  "
$ws.Range("D48").Value = "'0.132"
$ws.Range("E48").Value = "'  +3.96%  "
$ws.Range("B49").Value = "'ApeXProtocol"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'3.09"
$ws.Range("E49").Value = "'  -0.69%  "
$ws.Range("B50").Value = "'THORChain"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'8.70"
$ws.Range("E50").Value = "'  +6.92%  "
$ws.Range("D51").Value = "'139.43"
$ws.Range("E51").Value = "'  +2.17%  "
